$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Date Last Updated column stores its dates as plain text (e.g. "02/15/2024"),
# so format the cells as Text first to stop Excel from auto-converting the
# replacement strings into date serial numbers.
$ws.Range("D10:D12").NumberFormat = "@"

# Rows 10-12 (Tasks 9, 11, 12) move from "In Development" to "Complete", which
# uses the same bold/colored font as the other "Complete" rows (e.g. B2). Copy
# that formatting over before writing the new status text.
$ws.Range("B2").Copy()
$ws.Range("B10:B12").PasteSpecial(-4122)  # xlPasteFormats

# Task 9: Register Backend - now Complete, note updated, date updated
$ws.Range("B10").Value = "Complete"
$ws.Range("C10").Value = "modal wasn’t closing after form misinput bug fixed, UI needs a few more polishing"
$ws.Range("D10").Value = "02/21/2024"

# Task 11: Logout Backend - now Complete, note updated, date updated
$ws.Range("B11").Value = "Complete"
$ws.Range("C11").Value = "fixed bug where flash mgs exists, UI needs a few more polishing"
$ws.Range("D11").Value = "02/21/2024"

# Task 12: Login Backend - now Complete, note updated, date updated
$ws.Range("B12").Value = "Complete"
$ws.Range("C12").Value = "fixed bug where flash mgs exists, UI needs a few more polishing"
$ws.Range("D12").Value = "02/21/2024"

# Task 13 now has a description
$ws.Range("A13").Value = "Task 13: Polish form sizes for different screens"

# Update the active selection to A14, matching the saved view state
$ws.Range("A14").Select()
